$wb = $excel.ActiveWorkbook

# Rename sheets (tabs) to reflect new task-order identifiers
$wb.Worksheets.Item(1).Name = "GNG_TO-16509960800968807"
$wb.Worksheets.Item(2).Name = "NB_TO-1650996081161948"
$wb.Worksheets.Item(3).Name = "RS_TO-1650996081161948"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509960812099452"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509960812739124"

# Sheet1 - GNG
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509960800648441.csv"
$ws1.Range("B3").Value = "GNG_stims-16509960800808794.csv"
$ws1.Range("B4").Value = "go_stims-16509960800808794.csv"
$ws1.Range("B5").Value = "GNG_stims-16509960800968807.csv"

# Sheet2 - NB
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16509960808819118.csv"
$ws2.Range("B3").Value = "ZB-match_3-16509960803859088.csv"
$ws2.Range("B4").Value = "ZB-match_7-16509960802328494.csv"
$ws2.Range("B5").Value = "TB-16509960809139476.csv"
$ws2.Range("B6").Value = "OB-16509960808099453.csv"
$ws2.Range("B7").Value = "TB-16509960811459513.csv"
$ws2.Range("B8").Value = "ZB-match_2-16509960803059454.csv"
$ws2.Range("B9").Value = "TB-16509960809699488.csv"
$ws2.Range("B10").Value = "OB-16509960808419454.csv"

# Sheet3 - RS
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# Sheet4 - TOL
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509960811779563.csv"
$ws4.Range("B3").Value = "ZM_stims-1650996081161948.csv"
$ws4.Range("B4").Value = "MM_stims-16509960811939442.csv"
$ws4.Range("B5").Value = "ZM_stims-16509960811779563.csv"
$ws4.Range("B6").Value = "MM_stims-16509960812099452.csv"
$ws4.Range("B7").Value = "ZM_stims-16509960811939442.csv"

# Sheet5 - vSAT
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16509960812579472.csv"
$ws5.Range("B3").Value = "SAT_stims-1650996081225915.csv"
$ws5.Range("B4").Value = "SAT_stims-16509960812099452.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509960812419171.csv"
